{"js": "// Replace the 25 three-digit-by-one-digit multiplication answers in the\n// document's table with a new set of problems, in document order.\nconst pairs = [\n  [\"616\u00d77=4312\", \"833\u00d74=3332\"],\n  [\"909\u00d78=7272\", \"763\u00d79=6867\"],\n  [\"215\u00d78=1720\", \"673\u00d78=5384\"],\n  [\"733\u00d79=6597\", \"275\u00d75=1375\"],\n  [\"875\u00d77=6125\", \"619\u00d77=4333\"],\n  [\"747\u00d78=5976\", \"916\u00d77=6412\"],\n  [\"334\u00d79=3006\", \"266\u00d72=532\"],\n  [\"695\u00d79=6255\", \"899\u00d75=4495\"],\n  [\"863\u00d72=1726\", \"885\u00d79=7965\"],\n  [\"195\u00d78=1560\", \"790\u00d73=2370\"],\n  [\"643\u00d74=2572\", \"467\u00d74=1868\"],\n  [\"129\u00d77=903\", \"422\u00d73=1266\"],\n  [\"387\u00d78=3096\", \"432\u00d77=3024\"],\n  [\"449\u00d73=1347\", \"592\u00d75=2960\"],\n  [\"146\u00d73=438\", \"366\u00d74=1464\"],\n  [\"319\u00d76=1914\", \"543\u00d74=2172\"],\n  [\"348\u00d76=2088\", \"458\u00d73=1374\"],\n  [\"584\u00d78=4672\", \"515\u00d76=3090\"],\n  [\"122\u00d74=488\", \"968\u00d78=7744\"],\n  [\"400\u00d72=800\", \"849\u00d78=6792\"],\n  [\"265\u00d78=2120\", \"304\u00d73=912\"],\n  [\"959\u00d75=4795\", \"472\u00d78=3776\"],\n  [\"339\u00d73=1017\", \"430\u00d75=2150\"],\n  [\"561\u00d78=4488\", \"721\u00d72=1442\"],\n  [\"410\u00d75=2050\", \"478\u00d76=2868\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n", "ps1": "# Replace the 25 three-digit-by-one-digit multiplication answers in the\n# document's table with a new set of problems, in document order.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"616\u00d77=4312\", \"833\u00d74=3332\"),\n  @(\"909\u00d78=7272\", \"763\u00d79=6867\"),\n  @(\"215\u00d78=1720\", \"673\u00d78=5384\"),\n  @(\"733\u00d79=6597\", \"275\u00d75=1375\"),\n  @(\"875\u00d77=6125\", \"619\u00d77=4333\"),\n  @(\"747\u00d78=5976\", \"916\u00d77=6412\"),\n  @(\"334\u00d79=3006\", \"266\u00d72=532\"),\n  @(\"695\u00d79=6255\", \"899\u00d75=4495\"),\n  @(\"863\u00d72=1726\", \"885\u00d79=7965\"),\n  @(\"195\u00d78=1560\", \"790\u00d73=2370\"),\n  @(\"643\u00d74=2572\", \"467\u00d74=1868\"),\n  @(\"129\u00d77=903\", \"422\u00d73=1266\"),\n  @(\"387\u00d78=3096\", \"432\u00d77=3024\"),\n  @(\"449\u00d73=1347\", \"592\u00d75=2960\"),\n  @(\"146\u00d73=438\", \"366\u00d74=1464\"),\n  @(\"319\u00d76=1914\", \"543\u00d74=2172\"),\n  @(\"348\u00d76=2088\", \"458\u00d73=1374\"),\n  @(\"584\u00d78=4672\", \"515\u00d76=3090\"),\n  @(\"122\u00d74=488\", \"968\u00d78=7744\"),\n  @(\"400\u00d72=800\", \"849\u00d78=6792\"),\n  @(\"265\u00d78=2120\", \"304\u00d73=912\"),\n  @(\"959\u00d75=4795\", \"472\u00d78=3776\"),\n  @(\"339\u00d73=1017\", \"430\u00d75=2150\"),\n  @(\"561\u00d78=4488\", \"721\u00d72=1442\"),\n  @(\"410\u00d75=2050\", \"478\u00d76=2868\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n  #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n  #   ReplaceWith, Replace)\n  # Wrap=1 (wdFindContinue), Replace=2 (wdReplaceOne) since each answer is unique.\n  $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
